# "begun grading HMWK 2" -- add a new "Homework 2" column (E) to the grade
# book and fill in the scores recorded so far (several students are still
# ungraded, matching the commit message).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column header
$ws.Range("E1").Value = "Homework 2"

# Scores entered so far, as fractions out of 25 (mirrors the existing
# Homework 1 column's "points/35" formula style). Rows left untouched
# (4, 5, 6, 10, 11, 12, 14) are students not yet graded for Homework 2.
$ws.Range("E2").Formula = "=26/25"
$ws.Range("E3").Value = " "
$ws.Range("E7").Value = 0
$ws.Range("E8").Formula = "=20/25"
$ws.Range("E9").Formula = "=22/25"
$ws.Range("E13").Formula = "=25/25"
$ws.Range("E15").Formula = "=17/25"

# Leave the selection where the grading work left off
$ws.Range("E14").Select()
